$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.620.77'
$ws.Range("E2").Value = '  -1.12%  '

$ws.Range("D3").Value = '3.504.02'
$ws.Range("E3").Value = '  +0.39%  '

$ws.Range("E4").Value = '  +0.07%  '

$ws.Range("D5").NumberFormat = '@'
$ws.Range("D5").Value = '582.29'
$ws.Range("E5").Value = '  -2.46%  '

$ws.Range("D6").NumberFormat = '@'
$ws.Range("D6").Value = '175.12'
$ws.Range("E6").Value = '  -2.73%  '

$ws.Range("E7").Value = '  +0.06%  '

$ws.Range("B8").Value = 'LidoStakedEther'
$ws.Range("C8").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("D8").Value = '3.506.82'
$ws.Range("E8").Value = '  +0.41%  '

$ws.Range("B9").Value = 'XRP'
$ws.Range("C9").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D9").NumberFormat = '@'
$ws.Range("D9").Value = '0.595'
$ws.Range("E9").Value = '  -1.92%  '

$ws.Range("D10").NumberFormat = '@'
$ws.Range("D10").Value = '0.134'
$ws.Range("E10").Value = '  -2.55%  '

$ws.Range("E11").Value = '  -1.87%  '

$ws.Range("D12").NumberFormat = '@'
$ws.Range("D12").Value = '0.423'
$ws.Range("E12").Value = '  -3.18%  '

$ws.Range("D13").Value = '4.113.40'
$ws.Range("E13").Value = '  +0.49%  '

$ws.Range("D14").NumberFormat = '@'
$ws.Range("D14").Value = '30.28'
$ws.Range("E14").Value = '  -6.09%  '

$ws.Range("D16").Value = '66.656.27'
$ws.Range("E16").Value = '  -1.05%  '

$ws.Range("E17").Value = '  -2.75%  '

$ws.Range("D18").Value = '3.524.49'
$ws.Range("E18").Value = '  +0.66%  '

$ws.Range("D19").NumberFormat = '@'
$ws.Range("D19").Value = '6.04'
$ws.Range("E19").Value = '  -4.00%  '

$ws.Range("D20").NumberFormat = '@'
$ws.Range("D20").Value = '13.97'
$ws.Range("E20").Value = '  -2.25%  '

$ws.Range("D21").NumberFormat = '@'
$ws.Range("D21").Value = '381.49'
$ws.Range("E21").Value = '  -2.19%  '

$ws.Range("D22").NumberFormat = '@'
$ws.Range("D22").Value = '7.85'
$ws.Range("E22").Value = '  -0.97%  '

$ws.Range("D23").NumberFormat = '@'
$ws.Range("D23").Value = '0.546'
$ws.Range("E23").Value = '  +0.58%  '

$ws.Range("E24").Value = '  +0.17%  '

$ws.Range("D25").NumberFormat = '@'
$ws.Range("D25").Value = '72.29'
$ws.Range("E25").Value = '  -2.27%  '

$ws.Range("E26").Value = '  +0.21%  '

$ws.Range("E27").Value = '  -1.38%  '

$ws.Range("E28").Value = '  -4.95%  '

$ws.Range("E29").Value = '  -1.89%  '

$ws.Range("E30").Value = '  -0.59%  '

$ws.Range("D31").NumberFormat = '@'
$ws.Range("D31").Value = '24.56'
$ws.Range("E31").Value = '  +4.26%  '

$ws.Range("D32").NumberFormat = '@'
$ws.Range("D32").Value = '5.86'
$ws.Range("E32").Value = '  -5.31%  '

$ws.Range("D33").NumberFormat = '@'
$ws.Range("D33").Value = '2.01'
$ws.Range("E33").Value = '  -2.31%  '

$ws.Range("E34").Value = '  -6.04%  '

$ws.Range("B35").Value = 'USDe'
$ws.Range("C35").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D35").NumberFormat = '@'
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  -0.02%  '

$ws.Range("B36").Value = 'Aptos'
$ws.Range("C36").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D36").NumberFormat = '@'
$ws.Range("D36").Value = '7.26'
$ws.Range("E36").Value = '  -1.73%  '

$ws.Range("E37").Value = '  -2.36%  '

$ws.Range("E38").Value = '  +12.79%  '

$ws.Range("D39").NumberFormat = '@'
$ws.Range("D39").Value = '160.64'
$ws.Range("E39").Value = '  -1.58%  '

$ws.Range("D40").NumberFormat = '@'
$ws.Range("D40").Value = '0.893'
$ws.Range("E40").Value = '  +2.71%  '

$ws.Range("D41").NumberFormat = '@'
$ws.Range("D41").Value = '1.78'
$ws.Range("E41").Value = '  -5.26%  '

$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D42").NumberFormat = '@'
$ws.Range("D42").Value = '6.52'
$ws.Range("E42").Value = '  -4.33%  '

$ws.Range("B43").Value = 'Filecoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D43").NumberFormat = '@'
$ws.Range("D43").Value = '4.50'
$ws.Range("E43").Value = '  -3.18%  '

$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").Value = '2.724.44'
$ws.Range("E44").Value = '  -4.42%  '

$ws.Range("B45").Value = 'dogwifhat'
$ws.Range("C45").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D45").NumberFormat = '@'
$ws.Range("D45").Value = '2.52'
$ws.Range("E45").Value = '  -10.25%  '

$ws.Range("D46").NumberFormat = '@'
$ws.Range("D46").Value = '0.0701'
$ws.Range("E46").Value = '  -2.91%  '

$ws.Range("B47").Value = 'OKB'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D47").NumberFormat = '@'
$ws.Range("D47").Value = '40.64'
$ws.Range("E47").Value = '  -2.54%  '

$ws.Range("B48").Value = 'InjectiveProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D48").NumberFormat = '@'
$ws.Range("D48").Value = '25.10'
$ws.Range("E48").Value = '  -6.16%  '

$ws.Range("E49").Value = '  -3.11%  '

$ws.Range("D50").NumberFormat = '@'
$ws.Range("D50").Value = '323.93'
$ws.Range("E50").Value = '  -2.98%  '

$ws.Range("E51").Value = '  -3.73%  '
